$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9202686041625157
$ws.Range("C2").Value = 0.2312851940432097
$ws.Range("D2").Value = 0.2216154181053156
$ws.Range("E2").Value = 0.173995006164489
$ws.Range("F2").Value = 1.158484340464028
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("I2").Value = 0.4743025055241112
$ws.Range("J2").Value = 0.1814400276531032
$ws.Range("N2").Value = 1.002392765118174
$ws.Range("O2").Value = 2.593177007097438
$ws.Range("B3").Value = 0.8229987116431516
$ws.Range("C3").Value = 0.2035736324707784
$ws.Range("D3").Value = 0.2168581183874068
$ws.Range("E3").Value = 0.1700559920835083
$ws.Range("F3").Value = 1.152771200007365
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("I3").Value = 0.4795121132207214
$ws.Range("J3").Value = 0.1771225843945246
$ws.Range("N3").Value = 1.006460705725544
$ws.Range("O3").Value = 2.588672475193846
$ws.Range("B4").Value = 0.7632777456283009
$ws.Range("C4").Value = 0.18652978349877
$ws.Range("D4").Value = 0.2140201439054294
$ws.Range("E4").Value = 0.1677234901444855
$ws.Range("F4").Value = 1.150013388387933
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("I4").Value = 0.4830496188176809
$ws.Range("J4").Value = 0.1745782619095451
$ws.Range("N4").Value = 1.0093724869295
$ws.Range("O4").Value = 2.587734933961769
$ws.Range("B5").Value = 0.7389432360139949
$ws.Range("C5").Value = 0.1795773785316612
$ws.Range("D5").Value = 0.2128846190596363
$ws.Range("E5").Value = 0.1667946530482389
$ws.Range("F5").Value = 1.149078023406872
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("I5").Value = 0.4845762209298741
$ws.Range("J5").Value = 0.1735682267089018
$ws.Range("N5").Value = 1.010663349112768
$ws.Range("O5").Value = 2.587811949816768
$ws.Range("B6").Value = 0.7349026898970124
$ws.Range("C6").Value = 0.1784225307833083
$ws.Range("D6").Value = 0.2126973355315442
$ws.Range("E6").Value = 0.1666417305556642
$ws.Range("F6").Value = 1.148934085324953
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("I6").Value = 0.4848348438394687
$ws.Range("J6").Value = 0.1734021293203085
$ws.Range("N6").Value = 1.010883999664891
$ws.Range("O6").Value = 2.58785244577922
$ws.Range("B7").Value = 0.7629495507727597
$ws.Range("C7").Value = 0.1864360483954215
$ws.Range("D7").Value = 0.2140047447752238
$ws.Range("E7").Value = 0.1677108757105259
$ws.Range("F7").Value = 1.150000010836536
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("I7").Value = 0.48306986298509
$ws.Range("J7").Value = 0.1745645317343403
$ws.Range("N7").Value = 1.009389473461695
$ws.Range("O7").Value = 2.58773411474138
$ws.Range("B8").Value = 0.886730149374273
$ws.Range("C8").Value = 0.2217364363929164
$ws.Range("D8").Value = 0.2199579240351284
$ws.Range("E8").Value = 0.172618989144631
$ws.Range("F8").Value = 1.156358676151896
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("I8").Value = 0.4760283728502941
$ws.Range("J8").Value = 0.1799292395575947
$ws.Range("N8").Value = 1.003709555393868
$ws.Range("O8").Value = 2.591243982878581
$ws.Range("B9").Value = 1.129435424059523
$ws.Range("C9").Value = 0.2907193810294757
$ws.Range("D9").Value = 0.2322874718310004
$ws.Range("E9").Value = 0.1829259138570762
$ws.Range("F9").Value = 1.174788671659798
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("I9").Value = 0.4649149978424667
$ws.Range("J9").Value = 0.1912965171573262
$ws.Range("N9").Value = 0.995849766821209
$ws.Range("O9").Value = 2.612667522433043
$ws.Range("B10").Value = 1.307680207747296
$ws.Range("C10").Value = 0.3412427635975064
$ws.Range("D10").Value = 0.2417421073336783
$ws.Range("E10").Value = 0.1909142272174336
$ws.Range("F10").Value = 1.191979143481916
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("I10").Value = 0.4584025110626868
$ws.Range("J10").Value = 0.2001673305943257
$ws.Range("N10").Value = 0.9920658500439572
$ws.Range("O10").Value = 2.637325798671242
$ws.Range("B11").Value = 1.388742744771719
$ws.Range("C11").Value = 0.3641906446379153
$ws.Range("D11").Value = 0.2461286260973594
$ws.Range("E11").Value = 0.1946386679148517
$ws.Range("F11").Value = 1.200595791623613
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("I11").Value = 0.4558007153939982
$ws.Range("J11").Value = 0.2043163326977293
$ws.Range("N11").Value = 0.9907751797901483
$ws.Range("O11").Value = 2.650492070846184
$ws.Range("B12").Value = 1.419434611331553
$ws.Range("C12").Value = 0.3728750208108522
$ws.Range("D12").Value = 0.2478019101151716
$ws.Range("E12").Value = 0.1960620151282697
$ws.Range("F12").Value = 1.203973474671258
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("I12").Value = 0.4548675317125586
$ws.Range("J12").Value = 0.2059038214718498
$ws.Range("N12").Value = 0.9903482356929914
$ws.Range("O12").Value = 2.655758922751204
$ws.Range("B13").Value = 1.412824805813614
$ws.Range("C13").Value = 0.3710049353032332
$ws.Range("D13").Value = 0.2474409972292761
$ws.Range("E13").Value = 0.1957548949851642
$ws.Range("F13").Value = 1.203240923791597
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("I13").Value = 0.4550661912411513
$ws.Range("J13").Value = 0.2055611999136033
$ws.Range("N13").Value = 0.9904374388548263
$ws.Range("O13").Value = 2.654612098670441
$ws.Range("B14").Value = 1.391267887626498
$ws.Range("C14").Value = 0.3649052257147787
$ws.Range("D14").Value = 0.246266044147518
$ws.Range("E14").Value = 0.1947555075551151
$ws.Range("F14").Value = 1.200871374868498
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("I14").Value = 0.4557228974114906
$ws.Range("J14").Value = 0.2044466085696257
$ws.Range("N14").Value = 0.9907388168662408
$ws.Range("O14").Value = 2.650919740319353
$ws.Range("B15").Value = 1.378062985703991
$ws.Range("C15").Value = 0.3611682511823915
$ws.Range("D15").Value = 0.2455479379304677
$ws.Range("E15").Value = 0.1941450435770733
$ws.Range("F15").Value = 1.199434906690925
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("I15").Value = 0.4561319340344205
$ws.Range("J15").Value = 0.2037660190387527
$ws.Range("N15").Value = 0.9909314648552936
$ws.Range("O15").Value = 2.648694691369883
$ws.Range("B16").Value = 1.302382002015122
$ws.Range("C16").Value = 0.3397423187672644
$ws.Range("D16").Value = 0.2414571512294117
$ws.Range("E16").Value = 0.1906726447961375
$ws.Range("F16").Value = 1.191432071721763
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("I16").Value = 0.4585798196457489
$ws.Range("J16").Value = 0.1998984710504175
$ws.Range("N16").Value = 0.9921588536034562
$ws.Range("O16").Value = 2.636504640428456
$ws.Range("B17").Value = 1.255947413905119
$ws.Range("C17").Value = 0.3265888439129014
$ws.Range("D17").Value = 0.2389694293152189
$ws.Range("E17").Value = 0.1885656000786184
$ws.Range("F17").Value = 1.186726761655933
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("I17").Value = 0.4601740488354764
$ws.Range("J17").Value = 0.1975549666486103
$ws.Range("N17").Value = 0.9930220282897437
$ws.Range("O17").Value = 2.629526214108353
$ws.Range("B18").Value = 1.22923750050893
$ws.Range("C18").Value = 0.3190199987059543
$ws.Range("D18").Value = 0.2375466164435522
$ws.Range("E18").Value = 0.1873622055931534
$ws.Range("F18").Value = 1.184095361538581
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("I18").Value = 0.4611249535400361
$ws.Range("J18").Value = 0.1962177390384028
$ws.Range("N18").Value = 0.993559045553809
$ws.Range("O18").Value = 2.625695818934162
$ws.Range("B19").Value = 1.220193694058025
$ws.Range("C19").Value = 0.3164567604062256
$ws.Range("D19").Value = 0.2370662635925669
$ws.Range("E19").Value = 0.1869562215739364
$ws.Range("F19").Value = 1.183217284922009
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("I19").Value = 0.4614527382894025
$ws.Range("J19").Value = 0.1957668130865926
$ws.Range("N19").Value = 0.9937478373309006
$ws.Range("O19").Value = 2.62443038981931
$ws.Range("B20").Value = 1.260890670250774
$ws.Range("C20").Value = 0.3279893996150633
$ws.Range("D20").Value = 0.2392334182016072
$ws.Range("E20").Value = 0.1887890168460871
$ws.Range("F20").Value = 1.187219889452749
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("I20").Value = 0.4600008257883097
$ws.Range("J20").Value = 0.1978033297983899
$ws.Range("N20").Value = 0.9929259468280947
$ws.Range("O20").Value = 2.630250091395141
$ws.Range("B21").Value = 1.397599818925642
$ws.Range("C21").Value = 0.3666970103298581
$ws.Range("D21").Value = 0.2466108259049946
$ws.Range("E21").Value = 0.1950486998015961
$ws.Range("F21").Value = 1.201564253133057
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("I21").Value = 0.4555285925505395
$ws.Range("J21").Value = 0.2047735471842742
$ws.Range("N21").Value = 0.9906486184260075
$ws.Range("O21").Value = 2.651996640992536
$ws.Range("B22").Value = 1.486918681282248
$ws.Range("C22").Value = 0.3919625013272139
$ws.Range("D22").Value = 0.2515034802018619
$ws.Range("E22").Value = 0.1992154198121696
$ws.Range("F22").Value = 1.211608037269755
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("I22").Value = 0.4529092587716335
$ws.Range("J22").Value = 0.2094243033143783
$ws.Range("N22").Value = 0.9895204448726389
$ws.Range("O22").Value = 2.667847928427165
$ws.Range("B23").Value = 1.439250657971343
$ws.Range("C23").Value = 0.3784809052398259
$ws.Range("D23").Value = 0.2488857061321568
$ws.Range("E23").Value = 0.1969846516558462
$ws.Range("F23").Value = 1.206186211408266
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("I23").Value = 0.4542794134584831
$ws.Range("J23").Value = 0.20693338166042
$ws.Range("N23").Value = 0.9900896537945982
$ws.Range("O23").Value = 2.659237598996754
$ws.Range("B24").Value = 1.258655867698508
$ws.Range("C24").Value = 0.327356229340495
$ws.Range("D24").Value = 0.2391140457328191
$ws.Range("E24").Value = 0.1886879852870464
$ws.Range("F24").Value = 1.186996716676319
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("I24").Value = 0.4600790329071494
$ws.Range("J24").Value = 0.1976910134091412
$ws.Range("N24").Value = 0.9929692582594782
$ws.Range("O24").Value = 2.629922260868454
$ws.Range("B25").Value = 1.063785789668316
$ws.Range("C25").Value = 0.2720846572629227
$ws.Range("D25").Value = 0.2288821935888024
$ws.Range("E25").Value = 0.180064593882463
$ws.Range("F25").Value = 1.169163066708492
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("I25").Value = 0.4676319549144665
$ws.Range("J25").Value = 0.1881303834649231
$ws.Range("N25").Value = 0.9976259095982911
$ws.Range("O25").Value = 2.605309621831253
